$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (not auto-converted number) for Price-column values that look numeric,
# matching the source workbook where these are stored as literal text strings.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D11", "D13", "D14", "D19", "D20", "D22", "D24", "D25", "D27", "D28", "D29", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.639.53"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "2.552.21"
$ws.Range("E3").Value = "  +5.06%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "571.27"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("D6").Value = "151.34"
$ws.Range("E6").Value = "  +9.10%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "2.546.80"
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("D11").Value = "5.76"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D13").Value = "0.360"
$ws.Range("E13").Value = "  +3.48%  "
$ws.Range("D14").Value = "28.51"
$ws.Range("E14").Value = "  +8.46%  "
$ws.Range("D15").Value = "3.007.61"
$ws.Range("E15").Value = "  +5.10%  "
$ws.Range("D16").Value = "63.538.29"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").Value = "2.523.01"
$ws.Range("E18").Value = "  +4.00%  "
$ws.Range("D19").Value = "11.72"
$ws.Range("E19").Value = "  +4.71%  "
$ws.Range("D20").Value = "341.38"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("E21").Value = "  +4.47%  "
$ws.Range("D22").Value = "6.87"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "66.27"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("E26").Value = "  +5.30%  "
$ws.Range("B27").Value = "SuiNetwork"
$ws.Range("C27").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D27").Value = "1.52"
$ws.Range("E27").Value = "  +13.36%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "8.65"
$ws.Range("E28").Value = "  +6.06%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +12.65%  "
$ws.Range("D31").Value = "0.0₃0835"
$ws.Range("E31").Value = "  +6.47%  "
$ws.Range("D32").Value = "1.89"
$ws.Range("E32").Value = "  +4.68%  "
$ws.Range("D33").Value = "178.33"
$ws.Range("E33").Value = "  +3.72%  "
$ws.Range("D34").Value = "1.58"
$ws.Range("E34").Value = "  +9.19%  "
$ws.Range("D35").Value = "424.05"
$ws.Range("E35").Value = "  +11.21%  "
$ws.Range("E36").Value = "  +2.58%  "
$ws.Range("D37").Value = "19.22"
$ws.Range("E37").Value = "  +3.60%  "
$ws.Range("D38").Value = "4.48"
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "1.79"
$ws.Range("E39").Value = "  +6.39%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "39.79"
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("D43").Value = "154.20"
$ws.Range("E43").Value = "  +6.32%  "
$ws.Range("D44").Value = "3.82"
$ws.Range("E44").Value = "  +4.58%  "
$ws.Range("D45").Value = "21.18"
$ws.Range("E45").Value = "  +2.38%  "
$ws.Range("D46").Value = "0.613"
$ws.Range("E46").Value = "  +3.78%  "
$ws.Range("D47").Value = "0.0531"
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("D48").Value = "0.0970"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").Value = "0.0240"
$ws.Range("E49").Value = "  +7.97%  "
$ws.Range("D50").Value = "18.74"
$ws.Range("E50").Value = "  +4.71%  "
$ws.Range("E51").Value = "  +7.97%  "
